$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: the "Price" column (D) stores numeric-looking values as literal
# text in this workbook (t="inlineStr" in the original OOXML). Assigning a
# plain numeric-looking string via COM Range.Value coerces the cell to a
# Number, which would change both the stored type and the serialized form.
# Prefixing with a leading apostrophe forces Excel to keep (or make) the
# cell Text, matching the source data's type. Resetting the style back to
# "Normal" afterwards drops the implicit quote-prefix formatting flag that
# Excel attaches, so the cell's style stays identical to its original
# (unstyled) state.
function Set-TextValue($ref, $val) {
    $ws.Range($ref).Value = "'" + $val
    $ws.Range($ref).Style = "Normal"
}

# Row 2 - BNB : price update only
Set-TextValue "D2" "245.46"

# Row 3 - OKB : price update only
Set-TextValue "D3" "22.42"

# Row 4 - HuobiToken : price update only
Set-TextValue "D4" "5.470"

# Row 5 - Cronos : price update only
Set-TextValue "D5" "0.05615"

# Row 6 : KuCoinToken -> GateToken (data shifted up from the old row 7 set
# of values; a new GateToken entry is now in this slot)
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D6" "3.383"
$ws.Range("E6").Value = "5GateTokenGT"

# Row 7 : MXToken -> KuCoinToken
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D7" "6.469"
$ws.Range("E7").Value = "6KuCoinTokenKCS"

# Row 8 : FTXToken -> MXToken
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.8054"
$ws.Range("E8").Value = "7MXTokenMX"

# Row 9 : WazirX -> FTXToken
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D9" "1.043"
$ws.Range("E9").Value = "8FTXTokenFTT"

# Row 10 : MandalaExchangeToken -> WazirX
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1422"
$ws.Range("E10").Value = "9WazirXWRX"

# Row 11 : LiechtensteinCryptoassetsExchange -> MandalaExchangeToken
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07277"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

# Row 12 : BitrueCoin -> LiechtensteinCryptoassetsExchange
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03181"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

# Row 13 : BitMartToken -> BitrueCoin
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.02938"
$ws.Range("E13").Value = "12BitrueCoinBTR"

# Row 14 : BitForexToken -> BitMartToken
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09268"
$ws.Range("E14").Value = "13BitMartTokenBMX"

# Row 15 : MCDex -> BitForexToken
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001661"
$ws.Range("E15").Value = "14BitForexTokenBF"

# Row 16 : CoinExToken -> MCDex
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.221"
$ws.Range("E16").Value = "15MCDexMCB"

# Row 17 : One -> CoinExToken
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04742"
$ws.Range("E17").Value = "16CoinExTokenCET"

# Row 18 : TigerCash -> One
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005824"
$ws.Range("E18").Value = "17OneONE"

# Row 19 : HotbitToken -> TigerCash
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D19" "0.006434"
$ws.Range("E19").Value = "18TigerCashTCH"

# Row 20 : BitKan -> HotbitToken
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D20" "0.005072"
$ws.Range("E20").Value = "19HotbitTokenHTB"

# Row 21 : NitroEx -> BitKan
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D21" "0.001053"
$ws.Range("E21").Value = "20BitKanKAN"

# Row 22 : LEO -> NitroEx
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D22" "0.0001503"
$ws.Range("E22").Value = "21NitroExNTX"

# Row 23 : GateToken -> LEO
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D23" "3.986"
$ws.Range("E23").Value = "22LEOLEO"

# Row 24 - BTSEToken : price update only (unchanged row position)
Set-TextValue "D24" "2.121"

# Row 26 - ProBitToken volume label change
$ws.Range("E26").Value = "25ProBitTokenPROB"

# Row 40 - IDEX
Set-TextValue "D40" "0.04153"

# Row 41 - KickToken
Set-TextValue "D41" "0.006896"

# Row 42 - CEJI volume label change
$ws.Range("E42").Value = "41CEJICEJIBestin24h"

# Row 43 - BKEXToken
Set-TextValue "D43" "0.1039"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.009572"

# Row 45 - CoinLion
Set-TextValue "D45" "0.00005655"

# Row 47 - CoinbaseStockToken
Set-TextValue "D47" "0.6816"

# Row 48 - BOLO
Set-TextValue "D48" "0.01493"
